# Correction génération données fichiers json :
# - décalage des dates (A) de 2015 vers 2017
# - mise à jour des valeurs (E) suite à correction des données manquantes
#   (remplacées par NULL dans la base -> valeurs recalculées)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Nouvelles valeurs de la colonne E (nombre) pour les lignes 3 à 63
$eValues = @{
    3 = 10
    4 = 20
    5 = 18
    6 = 18
    7 = 6
    8 = 6
    9 = 16
    10 = 17
    11 = 13
    12 = 11
    13 = 7
    14 = 15
    15 = 20
    16 = 10
    17 = 13
    18 = 8
    19 = 18
    20 = 15
    21 = 11
    22 = 5
    23 = 20
    24 = 14
    25 = 17
    26 = 18
    27 = 19
    28 = 11
    29 = 16
    30 = 19
    31 = 6
    32 = 19
    33 = 16
    34 = 14
    35 = 7
    36 = 11
    37 = 6
    38 = 13
    39 = 8
    40 = 13
    41 = 20
    42 = 10
    43 = 10
    44 = 20
    45 = 10
    46 = 18
    47 = 16
    48 = 12
    49 = 17
    50 = 10
    51 = 18
    52 = 19
    53 = 18
    54 = 16
    55 = 19
    56 = 6
    57 = 17
    58 = 10
    59 = 11
    60 = 16
    61 = 16
    62 = 18
    63 = 15
}

for ($row = 3; $row -le 63; $row++) {
    # Colonne A : décalage de l'année de formation 2015 -> 2017 (+20000)
    $aCell = $ws.Cells.Item($row, 1)
    $aCell.Value2 = $aCell.Value2 + 20000

    # Colonne E : nouvelle valeur
    $ws.Cells.Item($row, 5).Value2 = $eValues[$row]
}
